# Updated cryptos list on Wed Jun 14 20:51:29 UTC 2023 with GitHub Actions
#
# Applies the new Price (col D) / Volume(1h) (col E) figures scraped for
# this run. Price values are plain text in the sheet (e.g. "1.669.80" uses
# dots as thousands separators and would otherwise be auto-coerced into a
# number by Excel), so each one is written with a leading apostrophe to
# force text entry and the cell style is immediately reset back to
# "Normal" so no stray number-format/quote-prefix style sticks to the
# cell. The Volume(1h) percentages already contain a "%" plus padding
# spaces and stay text on their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($cellRef, $text)
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

function Set-Volume {
    param($cellRef, $text)
    $ws.Range($cellRef).Value = $text
}

# Row 2 - Bitcoin
Set-PriceText "D2" "25.228.85"
Set-Volume    "E2" "  -2.61%  "

# Row 3 - Ethereum
Set-PriceText "D3" "1.661.09"
Set-Volume    "E3" "  -4.63%  "

# Row 4 - TetherUSD
Set-PriceText "D4" "1.009"
Set-Volume    "E4" "  +1.03%  "

# Row 5 - BNB
Set-PriceText "D5" "237.07"
Set-Volume    "E5" "  -1.67%  "

# Row 6 - USDC (price unchanged)
Set-Volume    "E6" "  +0.77%  "

# Row 7 - XRP
Set-PriceText "D7" "0.4757"
Set-Volume    "E7" "  -8.98%  "

# Row 8 - Cardano
Set-PriceText "D8" "0.2612"
Set-Volume    "E8" "  -5.14%  "

# Row 9 - Dogecoin
Set-PriceText "D9" "0.05950"
Set-Volume    "E9" "  -3.59%  "

# Row 10 - TRON
Set-PriceText "D10" "0.07115"
Set-Volume    "E10" "  -1.07%  "

# Row 11 - WrappedEther
Set-PriceText "D11" "1.691.07"
Set-Volume    "E11" "  -2.66%  "

# Row 12 - Polygon
Set-PriceText "D12" "0.6206"
Set-Volume    "E12" "  -3.87%  "

# Row 13 - Solana
Set-PriceText "D13" "14.37"
Set-Volume    "E13" "  -4.31%  "

# Row 14 - Polkadot
Set-PriceText "D14" "4.624"
Set-Volume    "E14" "  +0.04%  "

# Row 15 - Litecoin
Set-PriceText "D15" "72.82"
Set-Volume    "E15" "  -6.20%  "

# Row 16 - Dai
Set-PriceText "D16" "1.005"
Set-Volume    "E16" "  +0.45%  "

# Row 17 - BinanceUSD
Set-PriceText "D17" "1.007"
Set-Volume    "E17" "  +0.79%  "

# Row 18 - WrappedBTC
Set-PriceText "D18" "25.337.93"
Set-Volume    "E18" "  -2.17%  "

# Row 19 - Avalanche
Set-PriceText "D19" "11.44"
Set-Volume    "E19" "  -2.61%  "

# Row 20 - ShibaInu
Set-PriceText "D20" "0.000006581"
Set-Volume    "E20" "  -2.99%  "

# Row 21 - WrappedliquidstakedEther2.0
Set-PriceText "D21" "1.909.59"
Set-Volume    "E21" "  -2.89%  "

# Row 22 - Uniswap
Set-PriceText "D22" "4.427"
Set-Volume    "E22" "  +3.10%  "

# Row 23 - Cosmos
Set-PriceText "D23" "8.574"
Set-Volume    "E23" "  -1.04%  "

# Row 24 - Chainlink
Set-PriceText "D24" "5.247"
Set-Volume    "E24" "  -0.83%  "

# Row 25 - Monero
Set-PriceText "D25" "133.25"
Set-Volume    "E25" "  -4.10%  "

# Row 26 - EthereumClassic
Set-PriceText "D26" "14.69"
Set-Volume    "E26" "  -3.41%  "

# Row 27 - Toncoin
Set-PriceText "D27" "1.374"
Set-Volume    "E27" "  -9.92%  "

# Row 28 - LidoDAOToken
Set-PriceText "D28" "1.707"
Set-Volume    "E28" "  -3.55%  "

# Row 29 - BitcoinCash
Set-PriceText "D29" "102.22"
Set-Volume    "E29" "  -3.76%  "

# Row 30 - InternetComputer(DFINITY)
Set-PriceText "D30" "3.828"
Set-Volume    "E30" "  -2.53%  "

# Row 31 - Stellar
Set-PriceText "D31" "0.07891"
Set-Volume    "E31" "  -5.26%  "

# Row 32 - Filecoin
Set-PriceText "D32" "3.520"
Set-Volume    "E32" "  -4.89%  "

# Row 33 - Hedera
Set-PriceText "D33" "0.04599"
Set-Volume    "E33" "  -0.74%  "

# Row 34 - HuobiToken
Set-PriceText "D34" "2.643"
Set-Volume    "E34" "  -0.03%  "

# Row 35 - ARBITRUM
Set-PriceText "D35" "0.9379"
Set-Volume    "E35" "  -5.45%  "

# Row 36 - ImmutableX
Set-PriceText "D36" "0.5804"
Set-Volume    "E36" "  -6.52%  "

# Row 37 - MXToken
Set-PriceText "D37" "2.643"
Set-Volume    "E37" "  -1.48%  "

# Row 38 - VeChain
Set-PriceText "D38" "0.01547"
Set-Volume    "E38" "  -3.68%  "

# Row 39 - PaxDollar
Set-PriceText "D39" "1.006"
Set-Volume    "E39" "  +0.55%  "

# Row 40 - TrustWalletToken
Set-PriceText "D40" "0.8369"
Set-Volume    "E40" "  +12.58%  "

# Row 41 - RenderToken
Set-PriceText "D41" "1.837"
Set-Volume    "E41" "  -5.07%  "

# Row 42 - Quant
Set-PriceText "D42" "98.51"
Set-Volume    "E42" "  +0.61%  "

# Row 43 - TheSandbox
Set-PriceText "D43" "0.3712"
Set-Volume    "E43" "  -3.88%  "

# Row 44 - FraxShare
Set-PriceText "D44" "4.880"
Set-Volume    "E44" "  -2.46%  "

# Row 45 - Algorand
Set-PriceText "D45" "0.1136"
Set-Volume    "E45" "  +0.17%  "

# Row 46 / 47 - Aptos and Cronos swap ranking order, values refreshed
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-PriceText "D46" "0.05184"
Set-Volume    "E46" "  -1.13%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-PriceText "D47" "6.026"
Set-Volume    "E47" "  -3.86%  "

# Row 48 - Aave
Set-PriceText "D48" "53.47"
Set-Volume    "E48" "  -2.55%  "

# Row 49 - Elrond
Set-PriceText "D49" "29.52"
Set-Volume    "E49" "  -3.19%  "

# Row 50 - TrueUSD
Set-PriceText "D50" "1.006"
Set-Volume    "E50" "  +0.49%  "

# Row 51 - EnergySwap
Set-PriceText "D51" "7.379"
Set-Volume    "E51" "  -2.93%  "
